$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

# Row 2
Set-TextValue $ws "D2" "26.073.10"
Set-TextValue $ws "E2" "  -0.38%  "

# Row 3
Set-TextValue $ws "D3" "1.652.20"
Set-TextValue $ws "E3" "  -0.47%  "

# Row 4
Set-TextValue $ws "E4" "  -0.10%  "

# Row 5
Set-TextValue $ws "D5" "217.33"
Set-TextValue $ws "E5" "  +0.24%  "

# Row 6
Set-TextValue $ws "D6" "0.5252"
Set-TextValue $ws "E6" "  +2.07%  "

# Row 7
Set-TextValue $ws "E7" "  -0.06%  "

# Row 8
Set-TextValue $ws "D8" "0.2594"
Set-TextValue $ws "E8" "  -1.71%  "

# Row 9
Set-TextValue $ws "D9" "0.06332"
Set-TextValue $ws "E9" "  +1.09%  "

# Row 10
Set-TextValue $ws "D10" "20.35"
Set-TextValue $ws "E10" "  -1.88%  "

# Row 11
Set-TextValue $ws "D11" "0.07792"
Set-TextValue $ws "E11" "  +0.72%  "

# Row 12
Set-TextValue $ws "B12" "WrappedEther"
Set-TextValue $ws "C12" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws "D12" "1.688.82"
Set-TextValue $ws "E12" "  +1.54%  "

# Row 13
Set-TextValue $ws "B13" "Polkadot"
Set-TextValue $ws "C13" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws "D13" "4.494"
Set-TextValue $ws "E13" "  +1.10%  "

# Row 14
Set-TextValue $ws "D14" "0.5473"
Set-TextValue $ws "E14" "  +0.72%  "

# Row 15
Set-TextValue $ws "E15" "  +1.61%  "

# Row 16
Set-TextValue $ws "D16" "65.43"
Set-TextValue $ws "E16" "  +0.83%  "

# Row 17
Set-TextValue $ws "D17" "26.112.58"
Set-TextValue $ws "E17" "  -0.32%  "

# Row 18
Set-TextValue $ws "D18" "1.003"
Set-TextValue $ws "E18" "  -0.09%  "

# Row 19
Set-TextValue $ws "D19" "4.571"
Set-TextValue $ws "E19" "  -1.22%  "

# Row 20
Set-TextValue $ws "D20" "190.83"
Set-TextValue $ws "E20" "  -0.93%  "

# Row 21
Set-TextValue $ws "E21" "  -0.12%  "

# Row 22
Set-TextValue $ws "D22" "6.019"
Set-TextValue $ws "E22" "  +0.09%  "

# Row 23
Set-TextValue $ws "E23" "  -0.08%  "

# Row 24
Set-TextValue $ws "D24" "142.21"
Set-TextValue $ws "E24" "  +1.60%  "

# Row 25
Set-TextValue $ws "D25" "0.1231"
Set-TextValue $ws "E25" "  +0.71%  "

# Row 26
Set-TextValue $ws "D26" "7.223"
Set-TextValue $ws "E26" "  -0.08%  "

# Row 27
Set-TextValue $ws "E27" "  -0.84%  "

# Row 28
Set-TextValue $ws "D28" "1.430"
Set-TextValue $ws "E28" "  -0.26%  "

# Row 29
Set-TextValue $ws "D29" "0.05848"
Set-TextValue $ws "E29" "  -1.54%  "

# Row 30
Set-TextValue $ws "D30" "1.272"
Set-TextValue $ws "E30" "  +0.06%  "

# Row 31
Set-TextValue $ws "D31" "3.526"
Set-TextValue $ws "E31" "  -1.06%  "

# Row 32
Set-TextValue $ws "D32" "3.255"
Set-TextValue $ws "E32" "  -0.05%  "

# Row 33
Set-TextValue $ws "E33" "  -1.33%  "

# Row 34
Set-TextValue $ws "D34" "0.9471"
Set-TextValue $ws "E34" "  -1.86%  "

# Row 35
Set-TextValue $ws "D35" "2.412"
Set-TextValue $ws "E35" "  -0.44%  "

# Row 36
Set-TextValue $ws "D36" "2.781"

# Row 37
Set-TextValue $ws "D37" "0.5716"
Set-TextValue $ws "E37" "  +1.19%  "

# Row 38
Set-TextValue $ws "E38" "  +1.33%  "

# Row 39
Set-TextValue $ws "D39" "5.764"
Set-TextValue $ws "E39" "  -3.29%  "

# Row 40
Set-TextValue $ws "D40" "0.8438"
Set-TextValue $ws "E40" "  -1.50%  "

# Row 41
Set-TextValue $ws "E41" "  -0.01%  "

# Row 42
Set-TextValue $ws "D42" "103.40"
Set-TextValue $ws "E42" "  +3.17%  "

# Row 43
Set-TextValue $ws "D43" "1.024.60"
Set-TextValue $ws "E43" "  +1.27%  "

# Row 44
Set-TextValue $ws "D44" "1.797.51"
Set-TextValue $ws "E44" "  -0.23%  "

# Row 45
Set-TextValue $ws "D45" "56.98"
Set-TextValue $ws "E45" "  +0.61%  "

# Row 46
Set-TextValue $ws "D46" "0.9998"
Set-TextValue $ws "E46" "  -0.75%  "

# Row 47
Set-TextValue $ws "D47" "0.4315"
Set-TextValue $ws "E47" "  +3.15%  "

# Row 48
Set-TextValue $ws "E48" "  -0.35%  "

# Row 49
Set-TextValue $ws "D49" "1.464"
Set-TextValue $ws "E49" "  +0.98%  "

# Row 50
Set-TextValue $ws "D50" "7.810"
Set-TextValue $ws "E50" "  -2.36%  "

# Row 51
Set-TextValue $ws "D51" "0.09640"
Set-TextValue $ws "E51" "  -0.63%  "

